$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Rule name in B11 ("R40") is being renamed to the literal text "1".
# A plain Value assignment of a numeric-looking string ("1") would be
# auto-converted to a number by Excel, which is not what we want here -
# the target cell must keep storing a text value ("1"), with its existing
# style untouched. We build the literal text via a formula (so Excel's
# input-parser never sees a bare numeric literal), then copy/paste the
# computed value over the target cell - this preserves the cell's
# existing number format/style exactly, unlike typing an apostrophe-
# prefixed string (which would mark the cell with a quote-prefix and
# allocate a new style).
$helper = $ws.Range("Z1")
$helper.Formula = "=""1"""
$helper.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$helper.Clear()
